$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = 9.23076923076923
$ws.Activate() | Out-Null
$ws.Range("S6").Select() | Out-Null
